# Vitae SharedResource.xlsx - "Error message popups added"
#
# Changes:
#  1. Fix German typo: Grossbuchstabenzeichen (ss instead of ss-sharp/eszett) on
#     the "de" sheet, row 99 column B (key: PasswordRequiresUpper).
#  2. Add two new resource rows to both the "en" and "de" sheets:
#       TheLogout        -> Logout / Abmeldung
#       ResetPasswordBy  -> Please reset your password by /
#                           Bitte setzen Sie ihr Passwort zurück in dem Sie
#  3. Make the "de" sheet the active tab/sheet, with A172:B172 selected on
#     both sheets (matching the newly appended last row).

$wb = $excel.ActiveWorkbook
$ws_en = $wb.Sheets.Item("en")
$ws_de = $wb.Sheets.Item("de")

# 1. Fix the German typo ('ß' -> 'ss') in the existing password-rule message.
$ws_de.Range("B99").Value = "Passwörter müssen mindestens ein Grossbuchstabenzeichen ('A'-'Z') haben."

# 2. Append the two new localisation rows (row 171 and 172) on both sheets.
$ws_en.Range("A171").Value = "TheLogout"
$ws_en.Range("B171").Value = "Logout"
$ws_en.Range("A172").Value = "ResetPasswordBy"
$ws_en.Range("B172").Value = "Please reset your password by"

$ws_de.Range("A171").Value = "TheLogout"
$ws_de.Range("B171").Value = "Abmeldung"
$ws_de.Range("A172").Value = "ResetPasswordBy"
$ws_de.Range("B172").Value = "Bitte setzen Sie ihr Passwort zurück in dem Sie"

# Match the existing formatting used by the rest of column A in this block
# (wrap text, vertically centred) on the newly added key cells by copying
# the format from the row above (A170), same as the rest of that column.
$ws_en.Range("A170").Copy()
$ws_en.Range("A171:A172").PasteSpecial(-4122)  # xlPasteFormats
$ws_de.Range("A170").Copy()
$ws_de.Range("A171:A172").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Update view state: "de" becomes the active/selected tab, and both
#    sheets end up with the newly added last row (A172:B172) selected.
$ws_en.Range("A172:B172").Select()
$ws_de.Activate()
$ws_de.Range("A172:B172").Select()
